$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column for rows 2-5 from 45175 (2023-09-06)
# to 45183 (2023-09-14), preserving existing date formatting.
foreach ($r in 2..5) {
    $ws.Range("C$r").Value = 45183
}
